$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "Category" (column B) values for rows 88-125
$categories = @{
    88  = "Government"
    89  = "NGO"
    90  = "Government"
    91  = "Government"
    92  = "Government"
    93  = "NGO"
    94  = "NGO"
    95  = "NGO"
    96  = "Donor"
    97  = "Social enterprise"
    99  = "Government"
    100 = "NGO"
    101 = "Government"
    102 = "Government"
    103 = "Government"
    104 = "Government"
    105 = "Government"
    106 = "Government"
    107 = "Government"
    108 = "Government"
    109 = "Government"
    110 = "Government"
    111 = "Government"
    112 = "Government"
    113 = "Government"
    114 = "Government"
    115 = "Government"
    116 = "NGO"
    117 = "NGO"
    118 = "Government"
    119 = "NGO"
    120 = "NGO"
    121 = "NGO"
    122 = "NGO"
    123 = "NGO"
    124 = "NGO"
    125 = "Government"
}

foreach ($row in $categories.Keys) {
    $ws.Cells.Item($row, 2).Value = $categories[$row]
}

# Re-apply an AutoFilter over the data range (matches the saved file's extent)
$ws.Range("A1:B125").AutoFilter(1) | Out-Null

# Make sure the sheet-scoped hidden _FilterDatabase defined name exists/is correct
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Sheet 1'!`$A`$1:`$B`$125")
$filterName.Visible = $false

# Add the new implementor row (after the filter so its range is not expanded)
$ws.Range("A126").Value = "Ministry of Local Government"
$ws.Range("B126").Value = "Government"

# Selection state to match the saved view
$ws.Range("B126").Select() | Out-Null
